$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to remain plain text so that number-like strings
    # (e.g. "1.00", "0.163") are not silently coerced into numeric values
    # by Excel's automatic type inference.
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "65.139.14"
Set-TextValue "E2" "  +4.50%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.489.73"
Set-TextValue "E3" "  +1.68%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.06%  "

# Row 5 - BNB
Set-TextValue "D5" "416.57"
Set-TextValue "E5" "  +0.75%  "

# Row 6 - Solana
Set-TextValue "D6" "130.91"
Set-TextValue "E6" "  +1.69%  "

# Row 7 - XRP
Set-TextValue "D7" "0.654"
Set-TextValue "E7" "  +4.62%  "

# Row 8 - USDC (only E changes)
Set-TextValue "E8" "  +0.05%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.778"
Set-TextValue "E9" "  +7.13%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.163"
Set-TextValue "E10" "  +15.67%  "

# Row 11 - Avalanche
Set-TextValue "D11" "43.28"
Set-TextValue "E11" "  +1.26%  "

# Row 12 - ShibaInu
Set-TextValue "D12" "0.0000266"
Set-TextValue "E12" "  +20.64%  "

# Row 13 - Polkadot
Set-TextValue "D13" "10.02"
Set-TextValue "E13" "  +8.67%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue "D14" "4.044.17"
Set-TextValue "E14" "  +1.86%  "

# Row 15 - TRON (only E changes)
Set-TextValue "E15" "  +0.08%  "

# Row 16 - Chainlink
Set-TextValue "D16" "20.37"
Set-TextValue "E16" "  -0.40%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "3.506.03"
Set-TextValue "E17" "  +2.93%  "

# Row 18 - Uniswap
Set-TextValue "D18" "12.91"
Set-TextValue "E18" "  +1.47%  "

# Row 19 - Polygon (only E changes)
Set-TextValue "E19" "  +2.13%  "

# Row 20 - WrappedBTC
Set-TextValue "D20" "65.083.89"
Set-TextValue "E20" "  +4.35%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "443.51"
Set-TextValue "E21" "  -6.84%  "

# Row 22 - Litecoin
Set-TextValue "D22" "89.43"
Set-TextValue "E22" "  -2.51%  "

# Row 23 - ImmutableX
Set-TextValue "D23" "3.23"
Set-TextValue "E23" "  -0.94%  "

# Row 24 - InternetComputer(DFINITY)
Set-TextValue "D24" "13.15"
Set-TextValue "E24" "  +0.26%  "

# Row 25 - PancakeSwap (only E changes)
Set-TextValue "E25" "  +1.12%  "

# Row 26 - Filecoin
Set-TextValue "D26" "9.91"
Set-TextValue "E26" "  +1.42%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "33.87"
Set-TextValue "E27" "  +1.46%  "

# Rows 28 & 29 - Toncoin and Cosmos swap ranking order
Set-TextValue "B28" "Cosmos"
Set-TextValue "C28" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D28" "12.47"
Set-TextValue "E28" "  +5.02%  "

Set-TextValue "B29" "Toncoin"
Set-TextValue "C29" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D29" "2.72"
Set-TextValue "E29" "  +3.62%  "

# Row 30 - RenderToken
Set-TextValue "D30" "7.38"
Set-TextValue "E30" "  -4.20%  "

# Row 31 - Hedera (only E changes)
Set-TextValue "E31" "  +4.98%  "

# Row 32 - Kaspa
Set-TextValue "D32" "0.162"
Set-TextValue "E32" "  -1.99%  "

# Row 33 - Dai
Set-TextValue "D33" "1.00"
Set-TextValue "E33" "  +0.03%  "

# Row 34 - InjectiveProtocol
Set-TextValue "D34" "39.35"
Set-TextValue "E34" "  -3.52%  "

# Row 35 - OKB
Set-TextValue "D35" "57.26"
Set-TextValue "E35" "  -1.14%  "

# Row 36 - VeChain
Set-TextValue "D36" "0.0504"
Set-TextValue "E36" "  +3.35%  "

# Row 37 - PEPE
Set-TextValue "D37" "0.0₃0723"
Set-TextValue "E37" "  +33.59%  "

# Row 38 - Stellar (only E changes)
Set-TextValue "E38" "  +8.76%  "

# Row 39 - FirstDigitalUSD
Set-TextValue "D39" "0.999"
Set-TextValue "E39" "  +0.00%  "

# Row 40 - WEMIXToken (only E changes)
Set-TextValue "E40" "  +4.95%  "

# Row 41 - Stacks (only E changes)
Set-TextValue "E41" "  -0.34%  "

# Row 42 - NEARProtocol
Set-TextValue "D42" "4.49"
Set-TextValue "E42" "  +4.36%  "

# Row 43 - Monero
Set-TextValue "D43" "146.85"
Set-TextValue "E43" "  +1.84%  "

# Row 44 - LidoDAOToken
Set-TextValue "D44" "3.26"
Set-TextValue "E44" "  -1.99%  "

# Row 45 - TheGraph (only E changes)
Set-TextValue "E45" "  -4.16%  "

# Row 46 - ARBITRUM (only E changes)
Set-TextValue "E46" "  -3.21%  "

# Row 47 - ThetaToken
Set-TextValue "D47" "2.33"
Set-TextValue "E47" "  -3.38%  "

# Row 48 - Cronos
Set-TextValue "D48" "0.146"
Set-TextValue "E48" "  +5.00%  "

# Row 49 - Celestia
Set-TextValue "D49" "15.69"
Set-TextValue "E49" "  -3.94%  "

# Row 50 - ApeXProtocol (only E changes)
Set-TextValue "E50" "  +8.70%  "

# Row 51 - EnergySwap
Set-TextValue "D51" "21.48"
Set-TextValue "E51" "  -3.72%  "
